# Append two new response rows (participant IDs 18 and 19) to the
# PostExperiment survey table, growing Table1 from A1:L18 to A1:L20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Row 19 (ID 18) -------------------------------------------------
$newRow19 = $lo.ListRows.Add()

# Copy formatting from the row above (row 18) so number formats /
# quote-prefix styling carry over to the freshly inserted row.
$ws.Range("A18:L18").Copy()
$ws.Range("A19:L19").PasteSpecial(-4122)

$r19 = $newRow19.Range
$r19.Cells.Item(1, 1).Value = 18
$r19.Cells.Item(1, 2).Value = 44614.5958796296
$r19.Cells.Item(1, 3).Value = 44614.5973726852
$r19.Cells.Item(1, 4).Value = "2249443W@student.gla.ac.uk"
$r19.Cells.Item(1, 5).Value = "Kieran Waugh"
$r19.Cells.Item(1, 6).Value = "'18"
$r19.Cells.Item(1, 7).Value = "Touch In The Air;Pinch Anywhere;Pinch on Circle;Dwell;"
$r19.Cells.Item(1, 8).Value = "it was time consuming"
$r19.Cells.Item(1, 9).Value = "It was easier"
$r19.Cells.Item(1, 10).Value = "Touch In The Air;Pinch Anywhere;Pinch on Circle;Dwell;"
$r19.Cells.Item(1, 11).Value = "because it was the easiest"
$r19.Cells.Item(1, 12).Value = "Because it was the hardest and most time consuming"

# --- Row 20 (ID 19) -------------------------------------------------
$newRow20 = $lo.ListRows.Add()

$ws.Range("A19:L19").Copy()
$ws.Range("A20:L20").PasteSpecial(-4122)

$r20 = $newRow20.Range
$r20.Cells.Item(1, 1).Value = 19
$r20.Cells.Item(1, 2).Value = 44617.5091898148
$r20.Cells.Item(1, 3).Value = 44617.5111574074
$r20.Cells.Item(1, 4).Value = "2249443W@student.gla.ac.uk"
$r20.Cells.Item(1, 5).Value = "Kieran Waugh"
$r20.Cells.Item(1, 6).Value = "'19"
$r20.Cells.Item(1, 7).Value = "Pinch Anywhere;Touch In The Air;Pinch on Circle;Dwell;"
$r20.Cells.Item(1, 8).Value = "I lost control quickly, activation takes time"
$r20.Cells.Item(1, 9).Value = "Quick, the fact your hand can be anywhere"
$r20.Cells.Item(1, 10).Value = "Pinch Anywhere;Dwell;Pinch on Circle;Touch In The Air;"
$r20.Cells.Item(1, 11).Value = "no specific hand position needed"
$r20.Cells.Item(1, 12).Value = "I needed to get really close to the screen"

Write-Output "Appended rows 19-20 to Table1 (now $($lo.Range.Address()))"
